# Add a "tournaments" column (boolean) to the games worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell I1, styled like the other header cells (A1:H1 use style index 1)
$ws.Range("I1").Value = "tournaments"
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Boolean values for rows 2-51 (games 1-50), in order
$tournaments = @($true,$true,$true,$true,$true,$true,$true,$false,$true,$true,$false,$true,$true,$true,$true,$true,$true,$true,$true,$true,$true,$true,$true,$true,$true,$false,$false,$true,$false,$true,$false,$true,$true,$true,$true,$true,$true,$true,$true,$true,$true,$true,$true,$true,$true,$true,$true,$true,$true,$true)

for ($i = 0; $i -lt $tournaments.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $tournaments[$i]
}
